# Remove the standalone "Picture 6" logo picture (id=7) from slide 1.
# (It sat at the end of the slide's shape tree, right after the last
# group shape, just before the slide's closing extLst.)
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$target = $null
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Id -eq 7 -and $shape.Name -eq "Picture 6") {
        $target = $shape
        break
    }
}

if ($target -eq $null) {
    # Fallback: the picture was the last top-level shape on the slide.
    $target = $s.Shapes.Item($s.Shapes.Count)
}

$target.Delete()
